# Update "想去人数" (column F) figures on the 展览 and 全部类型 sheets
# to reflect newly scraped counts (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Row -> new F-column value for sheet "展览"
$updatesExhibition = @{
    2  = 8448
    3  = 8091
    9  = 141
    10 = 188
    13 = 183
    14 = 2262
    16 = 70
    19 = 142
    20 = 89
}

# Row -> new F-column value for sheet "全部类型" (F14 differs slightly from 展览)
$updatesAll = @{
    2  = 8448
    3  = 8091
    9  = 141
    10 = 188
    13 = 183
    14 = 2264
    16 = 70
    19 = 142
    20 = 89
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updatesExhibition.Keys) {
    $ws1.Cells.Item($row, 6).Value = $updatesExhibition[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updatesAll.Keys) {
    $ws4.Cells.Item($row, 6).Value = $updatesAll[$row]
}
